$d = $word.ActiveDocument

$replacements = @(
    @("409÷9=", "437÷8="),
    @("752÷6=", "933÷2="),
    @("286÷7=", "390÷7="),
    @("229÷8=", "235÷5="),
    @("913÷7=", "290÷7="),
    @("903÷2=", "569÷5="),
    @("703÷4=", "154÷3="),
    @("248÷5=", "419÷9="),
    @("522÷4=", "554÷9="),
    @("678÷7=", "441÷6="),
    @("965÷8=", "776÷5="),
    @("852÷5=", "351÷5="),
    @("813÷2=", "210÷6="),
    @("478÷6=", "914÷7="),
    @("515÷4=", "282÷7="),
    @("907÷2=", "712÷8="),
    @("359÷9=", "291÷2="),
    @("587÷9=", "804÷6="),
    @("633÷8=", "747÷8="),
    @("672÷6=", "816÷7="),
    @("160÷9=", "419÷5="),
    @("238÷7=", "318÷6="),
    @("166÷4=", "881÷2="),
    @("455÷2=", "797÷2="),
    @("862÷3=", "796÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
